# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets.
# F4: 28 -> 29
# F5: 131 -> 133

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 29
    $ws.Range("F5").Value = 133
}
